# CryCompanywiseStockReport_1.xlsx - stock-count correction pass.
#
# The source report has no formulas (every cell is a literal <v>), so the
# authoring tool that produced this edit re-keyed a handful of item rows'
# "Closing Qty" (column F) -- and, for two rows that got re-sorted, their
# whole data tuple (B/D/E/F) -- then rippled the change through the
# per-vendor "Sub Total:" rows and the final "Sub Total:"/"Grand Total:"
# rows. Column G ("Value") is always Rate(D) * Qty(F).
#
# We reproduce that by (1) writing the new leaf values, (2) recomputing G
# for every touched row, then (3) walking the sheet once to recompute every
# "Sub Total:" row as the sum of the G values in its block, and finally the
# trailing "Sub Total:"/"Grand Total:" rows as the sum of all the per-vendor
# subtotals above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Leaf-level edits -------------------------------------------------
# Row=line item row. Only the keys present are overwritten; everything
# else on the row (item code/name, MRP, etc.) is untouched.
# Rows 391/392 swapped their entire B/D/E/F tuple (two SKUs that share a
# name got re-ordered), so those two list every column explicitly.
$itemEdits = @(
    @{Row=122; F=352},
    @{Row=123; F=89},
    @{Row=135; F=13},
    @{Row=140; F=51},
    @{Row=299; F=34},
    @{Row=337; F=85},
    @{Row=339; F=23},
    @{Row=340; F=8},
    @{Row=348; F=83},
    @{Row=365; F=12},
    @{Row=391; B=57077; D=93.08; E=111.2; F=1},
    @{Row=392; B=61610; D=102.71; E=122.71; F=385},
    @{Row=404; F=4},
    @{Row=409; F=185},
    @{Row=412; F=2},
    @{Row=415; F=174},
    @{Row=424; F=384},
    @{Row=425; F=271},
    @{Row=426; F=226},
    @{Row=432; F=8},
    @{Row=438; F=264},
    @{Row=440; F=83},
    @{Row=444; F=67},
    @{Row=445; F=121},
    @{Row=446; F=266},
    @{Row=452; F=287},
    @{Row=460; F=157},
    @{Row=467; F=27},
    @{Row=473; F=33},
    @{Row=482; F=41},
    @{Row=495; F=6},
    @{Row=496; F=121},
    @{Row=528; F=69},
    @{Row=554; F=95},
    @{Row=591; F=703},
    @{Row=594; F=348},
    @{Row=599; F=517},
    @{Row=602; F=1146},
    @{Row=606; F=401},
    @{Row=638; F=9},
    @{Row=659; F=13},
    @{Row=728; F=122},
    @{Row=729; F=148},
    @{Row=746; F=40},
    @{Row=756; F=63},
    @{Row=757; F=146},
    @{Row=758; F=155},
    @{Row=769; F=170},
    @{Row=770; F=135},
    @{Row=772; F=175},
    @{Row=784; F=6},
    @{Row=798; F=29},
    @{Row=850; F=19},
    @{Row=851; F=4},
    @{Row=854; F=14},
    @{Row=867; F=40},
    @{Row=878; F=2},
    @{Row=881; F=30},
    @{Row=892; F=1},
    @{Row=901; F=94},
    @{Row=902; F=142},
    @{Row=903; F=266},
    @{Row=924; F=312},
    @{Row=927; F=198},
    @{Row=941; F=66},
    @{Row=945; F=230},
    @{Row=993; F=34},
    @{Row=996; F=14}
)

$colIndex = @{ B = 2; D = 4; E = 5; F = 6 }

foreach ($edit in $itemEdits) {
    $r = $edit.Row

    foreach ($col in @('B','D','E','F')) {
        if ($edit.ContainsKey($col)) {
            $ws.Cells.Item($r, $colIndex[$col]).Value2 = $edit[$col]
        }
    }

    # Value (G) = Rate (D) * Qty (F), recomputed from whatever is now on
    # the row (handles both the simple Qty-only edits and the 391/392
    # full-row swap the same way).
    $rate = $ws.Cells.Item($r, 4).Value2
    $qty = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 7).Value2 = $rate * $qty
}

# --- 2) Ripple into "Sub Total:" rows ------------------------------------
# Walk the used range once. Column A holds either:
#   - a numeric S.No for a line-item row (sum its G into the running block)
#   - the literal "Sub Total:" closing a vendor block (or the final overall
#     subtotal, which closes a run of "Sub Total:" rows instead)
#   - anything else (vendor name header, blank, "Grand Total:") which does
#     not participate in the running sum.
$lastRow = $ws.UsedRange.Rows.Count

$blockSum = 0.0           # running sum of G for the current line-item block
$blockItemCount = 0        # how many line-item rows fed $blockSum
$subtotalSum = 0.0         # running sum of B for "Sub Total:" rows seen so far
$subtotalCount = 0         # how many "Sub Total:" rows fed $subtotalSum
$lastSubtotalValue = 0.0   # value written into the most recent "Sub Total:" row

for ($r = 1; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2

    if ($a -eq "Sub Total:") {
        if ($blockItemCount -gt 0) {
            # Closes a block of line items.
            $newValue = $blockSum
        } else {
            # Closes a run of other "Sub Total:" rows (the final overall
            # subtotal) -- only happens once, at the very end.
            $newValue = $subtotalSum
        }
        $ws.Cells.Item($r, 2).Value2 = $newValue

        $lastSubtotalValue = $newValue
        $subtotalSum += $newValue
        $subtotalCount += 1
        $blockSum = 0.0
        $blockItemCount = 0
    }
    elseif ($a -is [double] -or $a -is [int]) {
        $blockSum += $ws.Cells.Item($r, 7).Value2
        $blockItemCount += 1
    }
    elseif ($a -eq "Grand Total:") {
        $ws.Cells.Item($r, 2).Value2 = $lastSubtotalValue
    }
    # else: vendor-name header row / blank row -- ignore.
}
